$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 7.5
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""

# Row 3
$ws.Range("D3").Value = 9.859999999999999
$ws.Range("E3").Value = 10.04
$ws.Range("F3").Value = 12.75
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""

# Row 4
$ws.Range("B4").Value = 12.5
$ws.Range("C4").Value = 10.14
$ws.Range("E4").Value = 9.98
$ws.Range("F4").Value = 10.33
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""

# Row 5
$ws.Range("C5").Value = 9.960000000000001
$ws.Range("D5").Value = 10.02
$ws.Range("F5").Value = 11.2
$ws.Range("G5").Value = 10.55
$ws.Range("H5").Value = 6.6
$ws.Range("I5").Value = ""

# Row 6
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = 7.25
$ws.Range("D6").Value = 9.67
$ws.Range("E6").Value = 8.800000000000001
$ws.Range("G6").Value = 11

# Row 7
$ws.Range("B7").Value = ""
$ws.Range("E7").Value = 9.449999999999999
$ws.Range("F7").Value = 9

# Row 8
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = 13.4
$ws.Range("I8").Value = 7.5

# Row 9
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("H9").Value = 12.5
$ws.Range("J9").Value = 15

# Row 10
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("I10").Value = 5
